# Trade #66 closed at 2026-02-16 21:35:06 - momentum DOWN +0.000%
#
# This script applies the set of trading-log updates described by the
# commit: a leadlag trade (#46, row 36) gets closed out with its exit
# numbers filled in, a new momentum trade (#66) is logged as freshly
# OPENed, that same closed leadlag trade is appended to the "All Trades"
# ledger, and the Summary / Comparison roll-up stats are refreshed to
# reflect the now one-trade-heavier, one-win-rate-lower totals.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Writes a string into a cell while forcing Text storage, even when
    # the string looks like a number/date/percentage (Excel would
    # otherwise silently reinterpret "65.2%" as 0.652 formatted as a
    # percentage, or "2026-02-16" as a date serial). NumberFormat is
    # restored afterwards via ClearFormats so no stray style survives.
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "Summary" - OVERALL and leadlag roll-up rows
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Cells.Item(2, 3).Value = 46
Set-TextValue $wsSummary.Cells.Item(2, 4) "65.2%"
Set-TextValue $wsSummary.Cells.Item(2, 5) "+10.9167%"
Set-TextValue $wsSummary.Cells.Item(2, 6) "+0.2373%"

$wsSummary.Cells.Item(3, 3).Value = 53
Set-TextValue $wsSummary.Cells.Item(3, 4) "39.6%"
Set-TextValue $wsSummary.Cells.Item(3, 5) "+6.7151%"
Set-TextValue $wsSummary.Cells.Item(3, 6) "+0.1267%"

# ---------------------------------------------------------------------
# Sheet "leadlag" - trade #46 (row 36) goes from OPEN to CLOSED
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Cells.Item(36, 7).Value = 68738.40286
$wsLeadlag.Cells.Item(36, 8).Value = "CLOSED"
$wsLeadlag.Cells.Item(36, 9).Value = -0.1413
$wsLeadlag.Cells.Item(36, 10).Value = -1.41
$wsLeadlag.Cells.Item(36, 13).Value = "time_exit_5min"
$wsLeadlag.Cells.Item(36, 14).Value = 5

# ---------------------------------------------------------------------
# Sheet "momentum" - append new trade #66 (row 14), still OPEN
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

$wsMomentum.Cells.Item(14, 1).Value = 66
Set-TextValue $wsMomentum.Cells.Item(14, 2) "2026-02-16"
Set-TextValue $wsMomentum.Cells.Item(14, 3) "21:35:06"
$wsMomentum.Cells.Item(14, 4).Value = "momentum"
$wsMomentum.Cells.Item(14, 5).Value = "DOWN"
$wsMomentum.Cells.Item(14, 6).Value = 68483.695
$wsMomentum.Cells.Item(14, 8).Value = "OPEN"
$wsMomentum.Cells.Item(14, 9).Value = 0
$wsMomentum.Cells.Item(14, 10).Value = 0
$wsMomentum.Cells.Item(14, 11).Value = 0.9
$wsMomentum.Cells.Item(14, 12).Value = "Downward momentum: -0.343% over 10 samples"
$wsMomentum.Cells.Item(14, 14).Value = 0

# ---------------------------------------------------------------------
# Sheet "All Trades" - append the now-closed leadlag trade #46 (row 47)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(47, 1).Value = 46
Set-TextValue $wsAll.Cells.Item(47, 2) "2026-02-16"
Set-TextValue $wsAll.Cells.Item(47, 3) "21:30:05"
$wsAll.Cells.Item(47, 4).Value = "leadlag"
$wsAll.Cells.Item(47, 5).Value = "DOWN"
$wsAll.Cells.Item(47, 6).Value = 68641.43
$wsAll.Cells.Item(47, 7).Value = 68738.40286
$wsAll.Cells.Item(47, 8).Value = "CLOSED"
$wsAll.Cells.Item(47, 9).Value = -0.1413
$wsAll.Cells.Item(47, 10).Value = -1.41
$wsAll.Cells.Item(47, 11).Value = 0.75
$wsAll.Cells.Item(47, 12).Value = "Binance leading with -0.163% move"
$wsAll.Cells.Item(47, 13).Value = "time_exit_5min"
$wsAll.Cells.Item(47, 14).Value = 5

# ---------------------------------------------------------------------
# Sheet "Comparison" - leadlag roll-up row
# ---------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

$wsComparison.Cells.Item(2, 2).Value = 53
Set-TextValue $wsComparison.Cells.Item(2, 3) "39.6%"
Set-TextValue $wsComparison.Cells.Item(2, 4) "2.45"
Set-TextValue $wsComparison.Cells.Item(2, 6) "-0.3306%"
Set-TextValue $wsComparison.Cells.Item(2, 7) "1.63"

Write-Output "Applied trade #66 / trade #46 close-out updates"
